$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as TEXT (preserve the original inline/shared string type
# instead of letting Excel auto-coerce numeric- or percentage-looking strings into
# numbers). We flip the cell to the "@" text format, write the value, then restore
# the "Normal" style so no stray formatting is left behind.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '305.07'
Set-TextValue $ws.Range("E2") '0.65%'
# Row 3
Set-TextValue $ws.Range("E3") '3.12%'
# Row 4
Set-TextValue $ws.Range("E4") '-2.16%'
# Row 5
Set-TextValue $ws.Range("D5") '0.07832'
# Row 6
Set-TextValue $ws.Range("E6") '-4.68%'
# Row 7
Set-TextValue $ws.Range("D7") '8.006'
Set-TextValue $ws.Range("E7") '-0.62%'
# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range("D8") '0.9281'
Set-TextValue $ws.Range("E8") '0.17%'
# Row 9
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range("D9") '0.09808'
Set-TextValue $ws.Range("E9") '-3.02%'
# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range("D10") '0.1874'
Set-TextValue $ws.Range("E10") '2.90%'
# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range("D11") '0.08601'
Set-TextValue $ws.Range("E11") '1.39%'
# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range("D12") '0.03727'
Set-TextValue $ws.Range("E12") '10.50%'
# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range("D13") '0.09995'
Set-TextValue $ws.Range("E13") '0.81%'
# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range("D14") '0.001484'
Set-TextValue $ws.Range("E14") '1.16%'
# Row 15
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range("D15") '0.005713'
Set-TextValue $ws.Range("E15") '-1.11%'
# Row 16
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D16") '3.469'
Set-TextValue $ws.Range("E16") '-0.05%'
# Row 17
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range("D17") '4.035'
Set-TextValue $ws.Range("E17") '1.48%'
# Row 18
Set-TextValue $ws.Range("D18") '2.352'
Set-TextValue $ws.Range("E18") '10.60%'
# Row 19
Set-TextValue $ws.Range("D19") '0.3411'
Set-TextValue $ws.Range("E19") '-0.56%'
# Row 20
Set-TextValue $ws.Range("D20") '0.1325'
Set-TextValue $ws.Range("E20") '0.04%'
# Row 21
Set-TextValue $ws.Range("D21") '4.760'
Set-TextValue $ws.Range("E21") '5.06%'
# Row 22
Set-TextValue $ws.Range("E22") '-0.78%'
# Row 23
Set-TextValue $ws.Range("D23") '0.04630'
Set-TextValue $ws.Range("E23") '0.03%'
# Row 24
Set-TextValue $ws.Range("D24") '0.001256'
Set-TextValue $ws.Range("E24") '3.44%'
# Row 25
Set-TextValue $ws.Range("D25") '0.004470'
Set-TextValue $ws.Range("E25") '0.13%'
# Row 26
Set-TextValue $ws.Range("D26") '0.0001400'
Set-TextValue $ws.Range("E26") '8.09%'
# Row 27
Set-TextValue $ws.Range("D27") '0.0002736'
Set-TextValue $ws.Range("E27") '-19.27%'
# Row 39
Set-TextValue $ws.Range("D39") '0.01806'
Set-TextValue $ws.Range("E39") '2.76%'
# Row 40
Set-TextValue $ws.Range("D40") '0.04757'
Set-TextValue $ws.Range("E40") '0.53%'
# Row 41
Set-TextValue $ws.Range("D41") '0.008030'
Set-TextValue $ws.Range("E41") '1.52%'
# Row 42
Set-TextValue $ws.Range("D42") '0.1403'
Set-TextValue $ws.Range("E42") '-0.92%'
# Row 43
Set-TextValue $ws.Range("D43") '0.007611'
Set-TextValue $ws.Range("E43") '-13.50%'
# Row 44
Set-TextValue $ws.Range("D44") '0.002109'
Set-TextValue $ws.Range("E44") '-8.00%'
# Row 45
Set-TextValue $ws.Range("D45") '0.01015'
Set-TextValue $ws.Range("E45") '10.87%'
# Row 46
Set-TextValue $ws.Range("D46") '0.00006328'
Set-TextValue $ws.Range("E46") '4.37%'
# Row 47
Set-TextValue $ws.Range("D47") '0.00000000755'
Set-TextValue $ws.Range("E47") '0.95%'
# Row 48
Set-TextValue $ws.Range("D48") '0.0005838'
Set-TextValue $ws.Range("E48") '0.65%'
# Row 49
Set-TextValue $ws.Range("D49") '35.04'
Set-TextValue $ws.Range("E49") '797.66%'
# Row 50
Set-TextValue $ws.Range("D50") '0.002707'
Set-TextValue $ws.Range("E50") '0.93%'
# Row 51
Set-TextValue $ws.Range("D51") '0.00002114'
Set-TextValue $ws.Range("E51") '0.95%'
